$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 430, pushing the existing
# rows 430-434 down to 433-437 (formatting/styles copied from row above
# by Excel's native Insert behaviour).
$ws.Rows("430:432").Insert()

# ---- New row 430 ----
$ws.Range("A430").Value = 4
$ws.Range("B430").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C430").Value = "Los Lagos"
$ws.Range("D430").Value = 44656
$ws.Range("E430").Value = 10
$ws.Range("F430").Value = "Fruta"
$ws.Range("G430").Value = 100108
$ws.Range("H430").Value = "Tropicales y subtropicales"
$ws.Range("I430").Value = 100108006
$ws.Range("J430").Value = "Plátano"
$ws.Range("K430").Value = "Barraganete"
$ws.Range("L430").Value = "Primera"
$ws.Range("M430").Value = 300
$ws.Range("N430").Value = 24000
$ws.Range("O430").Value = 25000
$ws.Range("P430").Value = 24500
$ws.Range("Q430").Value = "$/caja 20 kilos"
$ws.Range("R430").Value = "Ecuador"
$ws.Range("S430").Value = 1225
$ws.Range("T430").Value = 20

# ---- New row 431 ----
$ws.Range("A431").Value = 4
$ws.Range("B431").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C431").Value = "Los Lagos"
$ws.Range("D431").Value = 44656
$ws.Range("E431").Value = 10
$ws.Range("F431").Value = "Fruta"
$ws.Range("G431").Value = 100108
$ws.Range("H431").Value = "Tropicales y subtropicales"
$ws.Range("I431").Value = 100108006
$ws.Range("J431").Value = "Plátano"
$ws.Range("K431").Value = "Sin especificar"
$ws.Range("L431").Value = "Pintón"
$ws.Range("M431").Value = 500
$ws.Range("N431").Value = 18000
$ws.Range("O431").Value = 18000
$ws.Range("P431").Value = 18000
$ws.Range("Q431").Value = "$/caja 20 kilos"
$ws.Range("R431").Value = "Ecuador"
$ws.Range("S431").Value = 900
$ws.Range("T431").Value = 20

# ---- New row 432 ----
$ws.Range("A432").Value = 4
$ws.Range("B432").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C432").Value = "Los Lagos"
$ws.Range("D432").Value = 44656
$ws.Range("E432").Value = 10
$ws.Range("F432").Value = "Fruta"
$ws.Range("G432").Value = 100108
$ws.Range("H432").Value = "Tropicales y subtropicales"
$ws.Range("I432").Value = 100108006
$ws.Range("J432").Value = "Plátano"
$ws.Range("K432").Value = "Sin especificar"
$ws.Range("L432").Value = "Primera Pintón"
$ws.Range("M432").Value = 1000
$ws.Range("N432").Value = 19000
$ws.Range("O432").Value = 20000
$ws.Range("P432").Value = 19500
$ws.Range("Q432").Value = "$/caja 20 kilos"
$ws.Range("R432").Value = "Ecuador"
$ws.Range("S432").Value = 975
$ws.Range("T432").Value = 20
